$p = $ppt.ActivePresentation
Write-Host ($ppt | Get-Member | Out-String)
